# Cox Exploration build up
# Insert a new "Predictor" column (B) that repeats each row's predictor
# label (previously only in column A), shifting the existing statistic
# columns (old B..I) one column to the right (new C..J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B; this shifts the old B:I data
# columns to C:J and pushes dimension/row spans out automatically.
$ws.Columns("B").Insert()

# The newly inserted column B picked up column A's bordered/bold header
# style for rows 2-5 (since A2:A5 had it) but nothing for row 1 (A1 was
# never populated). Normalize data rows 2-5 to the default "Normal"
# style so they have no explicit style, matching the rest of the value
# columns.
$ws.Range("B2:B5").Style = "Normal"

# Populate the new header cell and give it the same bordered/bold header
# formatting as the other header cells (copy from C1, the old header
# start).
$ws.Range("B1").Value = "Predictor"
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)

# Populate the new column B data cells with the same predictor label
# already present in column A for that row (read back via .Text, which
# resolves to the displayed string, rather than .Value).
$ws.Range("B2").Value = $ws.Range("A2").Text
$ws.Range("B3").Value = $ws.Range("A3").Text
$ws.Range("B4").Value = $ws.Range("A4").Text
$ws.Range("B5").Value = $ws.Range("A5").Text
